$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest scraped values.
# Numeric-looking Price values are written with a leading apostrophe (forcing
# text, like typing them directly into Excel) so they keep their exact decimal
# formatting instead of being parsed into floating point numbers; the cell
# style is then reset to Normal so no stray numbering format sticks around.

$ws.Range("D2").Value = "63.800.10"
$ws.Range("D3").Value = "2.737.26"
$ws.Range("E3").Value = "  -0.81%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'569.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.64%  "
$ws.Range("D6").Value = "'155.55"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.66%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("E8").Value = "  -2.07%  "
$ws.Range("E9").Value = "  -3.65%  "
$ws.Range("E10").Value = "  -0.46%  "
$ws.Range("E11").Value = "  -2.69%  "
$ws.Range("D12").Value = "'5.38"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -19.87%  "
$ws.Range("D13").Value = "3.223.06"
$ws.Range("E13").Value = "  +0.15%  "
$ws.Range("E14").Value = "  -0.40%  "
$ws.Range("D15").Value = "63.502.90"
$ws.Range("E15").Value = "  -0.93%  "
$ws.Range("E16").Value = "  -3.55%  "
$ws.Range("D17").Value = "2.739.67"
$ws.Range("E17").Value = "  -0.44%  "
$ws.Range("D18").Value = "'11.98"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.33%  "
$ws.Range("E19").Value = "  -2.15%  "
$ws.Range("D20").Value = "'351.78"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.54%  "
$ws.Range("D21").Value = "'6.69"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.33%  "
$ws.Range("D22").Value = "'0.995"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.35%  "
$ws.Range("D23").Value = "'0.532"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("D24").Value = "'64.56"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.79%  "
$ws.Range("D25").Value = "'0.168"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.64%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").Value = "'8.30"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.65%  "
$ws.Range("E28").Value = "  -4.85%  "
$ws.Range("E29").Value = "  -4.20%  "
$ws.Range("D30").Value = "'6.85"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.80%  "
$ws.Range("E31").Value = "  -2.10%  "
$ws.Range("E32").Value = "  -7.69%  "
$ws.Range("D33").Value = "'19.95"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.39%  "
$ws.Range("E34").Value = "  +0.18%  "
$ws.Range("E35").Value = "  -0.28%  "
$ws.Range("E36").Value = "  -1.13%  "
$ws.Range("E37").Value = "  -3.45%  "
$ws.Range("D38").Value = "'0.963"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.86%  "
$ws.Range("E39").Value = "  +8.61%  "
$ws.Range("D40").Value = "'4.08"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.67%  "
$ws.Range("E41").Value = "  -6.90%  "
$ws.Range("D42").Value = "'38.62"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.08%  "
$ws.Range("D43").Value = "'21.06"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.54%  "
$ws.Range("E44").Value = "  -1.94%  "
$ws.Range("E45").Value = "  -3.30%  "
$ws.Range("D46").Value = "'134.47"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.96%  "
$ws.Range("E47").Value = "  -4.20%  "
$ws.Range("E48").Value = "  -1.84%  "
$ws.Range("D49").Value = "'0.0999"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.25%  "
$ws.Range("E50").Value = "  -0.03%  "
$ws.Range("D51").Value = "'11.05"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.57%  "
